$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1368.2084
$ws.Range("I17").Value = 1088.5555
$ws.Range("J17").Value = 1536
$ws.Range("K17").Value = 3265.6665
$ws.Range("L17").Value = 4608
$ws.Range("M17").Value = -3097.6665
$ws.Range("N17").Value = -4944
$ws.Range("H103").Value = 6149.8667
$ws.Range("I103").Value = 580.8
$ws.Range("J103").Value = 17288
$ws.Range("K103").Value = 1742.4
$ws.Range("L103").Value = 51864
$ws.Range("M103").Value = -1156.4
$ws.Range("N103").Value = -53036
$ws.Range("H132").Value = 15938171
$ws.Range("I132").Value = 17929896
$ws.Range("K132").Value = 53789688
$ws.Range("M132").Value = -53787158
$ws.Range("H141").Value = 1812.6364
$ws.Range("I141").Value = 1722.4062
$ws.Range("J141").Value = 4700
$ws.Range("K141").Value = 5167.2186
$ws.Range("L141").Value = 14100
$ws.Range("M141").Value = 12.78139999999985
$ws.Range("N141").Value = -24460

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1024.2885
$ws.Range("I61").Value = 831.5217
$ws.Range("J61").Value = 2502.1667
$ws.Range("K61").Value = 831.5217
$ws.Range("L61").Value = 2502.1667
$ws.Range("M61").Value = -619.5217
$ws.Range("N61").Value = -2926.1667
$ws.Range("H74").Value = 1801.9445
$ws.Range("I74").Value = 1371.8966
$ws.Range("J74").Value = 3583.5715
$ws.Range("K74").Value = 1371.8966
$ws.Range("L74").Value = 3583.5715
$ws.Range("M74").Value = -497.8966
$ws.Range("N74").Value = -5331.5715
$ws.Range("H76").Value = 34000
$ws.Range("J76").Value = 34000
$ws.Range("L76").Value = 34000
$ws.Range("N76").Value = -34676
$ws.Range("H77").Value = 1801.9445
$ws.Range("I77").Value = 1371.8966
$ws.Range("J77").Value = 3583.5715
$ws.Range("K77").Value = 6859.483
$ws.Range("L77").Value = 17917.8575
$ws.Range("M77").Value = -2491.483
$ws.Range("N77").Value = -26653.8575
$ws.Range("H79").Value = 34000
$ws.Range("J79").Value = 34000
$ws.Range("L79").Value = 34000
$ws.Range("N79").Value = -36340
$ws.Range("H122").Value = 2316.4666
$ws.Range("I122").Value = 1442.0769
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 4326.2307
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -1876.2307
$ws.Range("N122").Value = -28900
$ws.Range("H134").Value = 30040.584
$ws.Range("J134").Value = 30040.584
$ws.Range("L134").Value = 30040.584
$ws.Range("N134").Value = -40180.584
$ws.Range("H136").Value = 1024.2885
$ws.Range("I136").Value = 831.5217
$ws.Range("J136").Value = 2502.1667
$ws.Range("K136").Value = 2494.5651
$ws.Range("L136").Value = 7506.500100000001
$ws.Range("M136").Value = 55.4349000000002
$ws.Range("N136").Value = -12606.5001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 30000
$ws.Range("J53").Value = 30000
$ws.Range("L53").Value = 30000
$ws.Range("N53").Value = -31148

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2532.8914
$ws.Range("I31").Value = 980.3939
$ws.Range("K31").Value = 980.3939
$ws.Range("M31").Value = -685.3939
$ws.Range("H34").Value = 2532.8914
$ws.Range("I34").Value = 980.3939
$ws.Range("K34").Value = 980.3939
$ws.Range("M34").Value = -778.3939
$ws.Range("H81").Value = 27000
$ws.Range("J81").Value = 27000
$ws.Range("L81").Value = 27000
$ws.Range("N81").Value = -28996
$ws.Range("H84").Value = 27000
$ws.Range("J84").Value = 27000
$ws.Range("L84").Value = 81000
$ws.Range("N84").Value = -90984
$ws.Range("H94").Value = 2144.0715
$ws.Range("J94").Value = 2047.4546
$ws.Range("L94").Value = 2047.4546
$ws.Range("N94").Value = -2949.4546
$ws.Range("H107").Value = 1063.4286
$ws.Range("I107").Value = 689
$ws.Range("K107").Value = 689
$ws.Range("M107").Value = 1231
$ws.Range("H122").Value = 2207.087
$ws.Range("I122").Value = 1649.421
$ws.Range("J122").Value = 4856
$ws.Range("K122").Value = 4948.263
$ws.Range("L122").Value = 14568
$ws.Range("M122").Value = -2498.263
$ws.Range("N122").Value = -19468
$ws.Range("H132").Value = 2879.9807
$ws.Range("I132").Value = 1899.5
$ws.Range("K132").Value = 5698.5
$ws.Range("M132").Value = -3168.5
$ws.Range("H137").Value = 50255
$ws.Range("J137").Value = 50255
$ws.Range("L137").Value = 50255
$ws.Range("N137").Value = -60455
$ws.Range("H139").Value = 45268
$ws.Range("J139").Value = 45268
$ws.Range("L139").Value = 45268
$ws.Range("N139").Value = -55548

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 54.22222
$ws.Range("I12").Value = 174.66667
$ws.Range("J12").Value = 30.133333
$ws.Range("K12").Value = 524.00001
$ws.Range("L12").Value = 90.39999900000001
$ws.Range("M12").Value = -351.00001
$ws.Range("N12").Value = -436.399999
$ws.Range("H92").Value = 2452.2
$ws.Range("I92").Value = 455
$ws.Range("J92").Value = 4734.7144
$ws.Range("K92").Value = 1365
$ws.Range("L92").Value = 14204.1432
$ws.Range("M92").Value = -117
$ws.Range("N92").Value = -16700.1432

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1620
$ws.Range("I113").Value = 1366.6666
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1366.6666
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 803.3334
$ws.Range("N113").Value = -6340
$ws.Range("H122").Value = 3246.3447
$ws.Range("I122").Value = 2845.8262
$ws.Range("J122").Value = 4781.6665
$ws.Range("K122").Value = 8537.4786
$ws.Range("L122").Value = 14344.9995
$ws.Range("M122").Value = -6087.4786
$ws.Range("N122").Value = -19244.9995
$ws.Range("H132").Value = 2584.8164
$ws.Range("I132").Value = 1760.5667
$ws.Range("K132").Value = 5281.7001
$ws.Range("M132").Value = -2751.7001

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 863.05884
$ws.Range("I107").Value = 834.4545000000001
$ws.Range("J107").Value = 915.5
$ws.Range("K107").Value = 2503.3635
$ws.Range("L107").Value = 2746.5
$ws.Range("M107").Value = -583.3635000000004
$ws.Range("N107").Value = -6586.5
$ws.Range("H122").Value = 4271.4
$ws.Range("I122").Value = 3037.3809
$ws.Range("J122").Value = 10750
$ws.Range("K122").Value = 9112.1427
$ws.Range("L122").Value = 32250
$ws.Range("M122").Value = -6662.1427
$ws.Range("N122").Value = -37150
$ws.Range("H132").Value = 7940894
$ws.Range("I132").Value = 8516.308000000001
$ws.Range("J132").Value = 11496787
$ws.Range("K132").Value = 25548.924
$ws.Range("L132").Value = 34490361
$ws.Range("M132").Value = -23018.924
$ws.Range("N132").Value = -34495421
